# Insert a new weekly Ají (Feria Lagunitas de Puerto Montt) record as row 111,
# pushing the existing rows 111..174 down to 112..175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 111 - shifts rows 111:174 down to 112:175
# and extends the used range to A1:R175.
$ws.Rows("111:111").Insert()

# Populate the new row 111 with the new data point.
$ws.Cells.Item(111, 1).Value  = 4
$ws.Cells.Item(111, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(111, 3).Value  = "Los Lagos"
$ws.Cells.Item(111, 4).Value  = 44529
$ws.Cells.Item(111, 5).Value  = 10
$ws.Cells.Item(111, 6).Value  = 100112021
$ws.Cells.Item(111, 7).Value  = "Ají"
$ws.Cells.Item(111, 8).Value  = "Inferno"
$ws.Cells.Item(111, 9).Value  = "Primera"
$ws.Cells.Item(111, 10).Value = 40
$ws.Cells.Item(111, 11).Value = 22000
$ws.Cells.Item(111, 12).Value = 22000
$ws.Cells.Item(111, 13).Value = 22000
$ws.Cells.Item(111, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(111, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(111, 16).Value = 1833
$ws.Cells.Item(111, 17).Value = 12
$ws.Cells.Item(111, 18).Value = "Hortaliza"
